# Applies "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across several bullet
# paragraphs in the resume, matching the target diff.

$d = $word.ActiveDocument

# wdColor value for RRGGBB hex 2C3E50 (COM colors are stored BGR-packed)
$metricColor = 5258796

function Highlight-Metrics {
    param($ParaIndex, $Metrics)

    $p = $d.Paragraphs.Item($ParaIndex)
    $paraStart = $p.Range.Start
    $text = $p.Range.Text

    $searchFrom = 0
    foreach ($metric in $Metrics) {
        $offset = $text.IndexOf($metric, $searchFrom)
        if ($offset -lt 0) {
            continue
        }
        $rangeStart = $paraStart + $offset
        $rangeEnd = $rangeStart + $metric.Length
        $r = $d.Range($rangeStart, $rangeEnd)
        $r.Font.Bold = 1
        $r.Font.Color = $metricColor
        $searchFrom = $offset + $metric.Length
    }
}

# "• Discovered systematic race coding errors ... from 23% to 64%"
Highlight-Metrics 10 @("23%", "64%")

# "• Utilized advanced sampling methods ... from ±4.2% to ±2.1% ... from 71% to 87% ..."
Highlight-Metrics 12 @("±4.2%", "±2.1%", "71%", "87%")

# "• Trigonometric algorithm ... by 73.5%, saving ... $4.7M and enabling ..."
Highlight-Metrics 13 @("73.5%", "$4.7M")

# "• Built real-time FEC analysis systems ... valued over $2 trillion"
Highlight-Metrics 14 @("$2")

# "• Modernized legacy ETL processes ... reducing processing time by 57%"
Highlight-Metrics 39 @("57%")

# "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
Highlight-Metrics 55 @("73.5%")

# "• $4.7M savings enabled nonprofit access"
Highlight-Metrics 56 @("$4.7M")

# "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
Highlight-Metrics 57 @("12,847")

Write-Output "done"
